$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing last row value (row 70, column B): 157952 -> 174090
$ws.Cells.Item(70, 2).Value = 174090

# Add new row 71 with a new quarterly period.
# Column A must be stored as a plain text shared string ("01-04-2021"),
# not an Excel-autodetected date. Build it via a formula (forces text
# result), then paste-special as values-only into the target cell so it
# lands as plain text without picking up any date number-format/style.
$helper = $ws.Cells.Item(71, 5)
$helper.Formula = "=""01-04-2021"""
$helper.Copy()
$target = $ws.Cells.Item(71, 1)
$target.PasteSpecial(-4163)
$helper.ClearContents()

$ws.Cells.Item(71, 2).Value = 176505
$ws.Cells.Item(71, 3).Value = 176211
